$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22; existing rows 22:84 shift down to 23:85
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new weekly price record
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C22").Value = "Metropolitana"
$ws.Range("D22").Value = 44742
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 100114007
$ws.Range("G22").Value = "Jengibre"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 13833
$ws.Range("N22").Value = "$/caja 13 kilos"
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 1064
$ws.Range("Q22").Value = 13
$ws.Range("R22").Value = "Hortaliza"
